$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style of A2 (bold/border/center) down to the new rows A4:A7
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A4:A7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 2: Equalized Odds_Race_(0.001, 0.001)
$ws.Range("A2").Value = 'Equalized Odds_Race_(0.001, 0.001)'
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 50
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 0.001
$ws.Range("F2").Value = 0.001
$ws.Range("G2").Value = 'Equalized Odds'
$ws.Range("H2").Value = 'Race'
$ws.Range("I2").Value = '[[-0.8242363929748535, 0.014159401248010444], [-0.8745492696762085, 0.04942125245736562], [-0.8483567237854004, 0.022308763410192875], [-0.8272618055343628, 0.015153982078193153], [-0.8723942041397095, 0.041401527976946725], [-0.8729329705238342, 0.04163277990565041], [-0.8464503288269043, 0.019675899769046053], [-0.8696174621582031, 0.025799807076898284], [-0.873927652835846, 0.048502003441224226], [-0.8731402158737183, 0.043939038281920376], [-0.8715653419494629, 0.02747662629584531], [-0.872269868850708, 0.031295781863708405], [-0.8646027445793152, 0.024152517929055446]]'
$ws.Range("J2").Value = '[[0.0720604168948716], [0.002836020560472824], [0.0559131138617306], [0.08765127631437424], [0.0057057463996129515], [0.002577534814659025], [0.030659470640097105], [0.0028352034648014795], [0.003095331604070015], [0.005586780177463361], [0.001933568534453362], [0.0035356390842225284], [0.004155429204476469]]'

# Row 3: Equalized Odds_Race_(0.001, 0.1)
$ws.Range("A3").Value = 'Equalized Odds_Race_(0.001, 0.1)'
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 50
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 0.001
$ws.Range("F3").Value = 0.1
$ws.Range("G3").Value = 'Equalized Odds'
$ws.Range("H3").Value = 'Race'
$ws.Range("I3").Value = '[[-0.8729329705238342, 0.03932652152938046], [-0.8735546469688416, 0.04492886483761692], [-0.8708608150482178, 0.02273396964944826], [-0.8726842999458313, 0.039115291122808635], [-0.8721041083335876, 0.034613394954087606], [-0.8681669235229492, 0.02031283862696687]]'
$ws.Range("J3").Value = '[[0.004001572844964948], [0.002443965282760008], [0.003986077195804353], [0.0018077342267215526], [0.0039166856703735506], [0.001969444193204859]]'

# Row 4: Equalized Odds_Race_(0.001, 10)
$ws.Range("A4").Value = 'Equalized Odds_Race_(0.001, 10)'
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 50
$ws.Range("D4").Value = 10
$ws.Range("E4").Value = 0.001
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = 'Equalized Odds'
$ws.Range("H4").Value = 'Race'
$ws.Range("I4").Value = '[[-0.8619089126586914, 0.02549898203282666], [-0.8486883044242859, 0.02446810800867916], [-0.8347631692886353, 0.021626022809334945], [-0.872725784778595, 0.043879475927618386], [-0.8736375570297241, 0.04493837338658924], [-0.8678354024887085, 0.03005042979478375], [-0.8674623966217041, 0.02795540182508562], [-0.8722284436225891, 0.035812945280075156], [-0.8711509108543396, 0.03206587388804802], [-0.8726428747177124, 0.03783364585991919], [-0.8707779049873352, 0.030407067342858418], [-0.8629035353660583, 0.025808813413777], [-0.874880850315094, 0.05666187857656621], [-0.8680840730667114, 0.030337494227490495], [-0.8720211982727051, 0.03211967754367664]]'
$ws.Range("J4").Value = '[[0.009924848699921704], [0.07503941706300225], [0.09086869075900025], [0.004273323343272307], [0.0023895398211110945], [0.007803082023880141], [0.0026085558453432882], [0.006096653228117211], [0.0028502520542225845], [0.004606693194367068], [0.0038426154255498304], [0.010921646514229267], [0.0026291574446831944], [0.0032930808050942997], [0.005800639842243884]]'

# Row 5: Equalized Odds_Race_(0.1, 0.001)
$ws.Range("A5").Value = 'Equalized Odds_Race_(0.1, 0.001)'
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 50
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = 0.1
$ws.Range("F5").Value = 0.001
$ws.Range("G5").Value = 'Equalized Odds'
$ws.Range("H5").Value = 'Race'
$ws.Range("I5").Value = '[[-0.7633967399597168, 0.00013064210595074792], [-0.8711094260215759, 0.030486149007198676], [-0.872477114200592, 0.036467662446109717], [-0.8698247075080872, 0.02036614007050188], [-0.8718968629837036, 0.03514071765237675], [-0.8382444381713867, 0.014053534938677725], [-0.8739690780639648, 0.04636592800307515], [-0.8725599646568298, 0.042753248567297716]]'
$ws.Range("J5").Value = '[[0.08947120568403436], [0.0026459700330397585], [0.0038954644340870328], [0.004816235113004338], [0.00408031450697115], [0.03952480124322402], [0.002516430643276443], [0.0045271594241939675]]'

# Row 6: Equalized Odds_Race_(0.1, 0.1)
$ws.Range("A6").Value = 'Equalized Odds_Race_(0.1, 0.1)'
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 50
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 0.1
$ws.Range("F6").Value = 0.1
$ws.Range("G6").Value = 'Equalized Odds'
$ws.Range("H6").Value = 'Race'
$ws.Range("I6").Value = '[[-0.8728086352348328, 0.04566792015457151], [-0.7522483468055725, 0.0], [-0.8678354024887085, 0.03050767716561137], [-0.8716482520103455, 0.04253250961175356], [-0.8700318932533264, 0.0363622983488706], [-0.8728915452957153, 0.051737521594701505], [-0.8711509108543396, 0.040335867395115264], [-0.8417671918869019, 0.016279205014326607], [-0.8359235525131226, 0.010279678901658134], [-0.8548219799995422, 0.021896313358115156], [-0.8480252027511597, 0.021207311846490235], [-0.8722284436225891, 0.04377310740619205], [-0.8682498335838318, 0.035663286076179745], [-0.8658460974693298, 0.0295083420609425], [-0.8568941950798035, 0.02469360123870948], [-0.8636495471000671, 0.028323054558413448]]'
$ws.Range("J6").Value = '[[0.002836020560472824], [0.09395666529574745], [0.005660977500076756], [0.004329113995310084], [0.00421014211672947], [0.003668268240272974], [0.0025752790440976845], [0.054385473424672974], [0.058182731619724815], [0.024883423091979016], [0.061338221654791414], [0.004347405255735161], [0.009440373292847835], [0.005123701552146722], [0.03253262230350703], [0.003970822547819292]]'

# Row 7: Equalized Odds_Race_(0.1, 10)
$ws.Range("A7").Value = 'Equalized Odds_Race_(0.1, 10)'
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 50
$ws.Range("D7").Value = 10
$ws.Range("E7").Value = 0.1
$ws.Range("F7").Value = 10
$ws.Range("G7").Value = 'Equalized Odds'
$ws.Range("H7").Value = 'Race'
$ws.Range("I7").Value = '[[-0.7522483468055725, 0.0], [-0.8494343161582947, 0.001698615223809652], [-0.8671723008155823, 0.021765169040070798], [-0.8713995814323425, 0.035994645615542775], [-0.8741348385810852, 0.050330490712441084], [-0.8708608150482178, 0.032723838692287927], [-0.7522483468055725, 0.0], [-0.850801944732666, 0.01444244798396832], [-0.8734303116798401, 0.04998912041201213]]'
$ws.Range("J7").Value = '[[0.09579316406203515], [0.028936015992545895], [0.003589614050798086], [0.0031183331037854603], [0.002759758366152665], [0.002694870587713256], [0.09548765190984926], [0.06725469668650612], [0.002166291379806841]]'
